# Auto-generated edit script applying numeric corrections to Leve profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 1894.2858
$ws.Cells.Item(15, 9).Value = 1894.2858
$ws.Cells.Item(15, 11).Value = 5682.857400000001
$ws.Cells.Item(15, 13).Value = -5513.857400000001

$ws.Cells.Item(17, 8).Value = 3423.946
$ws.Cells.Item(17, 10).Value = 3491.2778
$ws.Cells.Item(17, 12).Value = 10473.8334
$ws.Cells.Item(17, 14).Value = -10809.8334

$ws.Cells.Item(20, 8).Value = 11000
$ws.Cells.Item(20, 10).Value = 14266.667
$ws.Cells.Item(20, 12).Value = 14266.667
$ws.Cells.Item(20, 14).Value = -14726.667

$ws.Cells.Item(26, 8).Value = 20400
$ws.Cells.Item(26, 9).Value = 0
$ws.Cells.Item(26, 10).Value = 20400
$ws.Cells.Item(26, 11).Value = 0
$ws.Cells.Item(26, 12).Value = 20400
$ws.Cells.Item(26, 14).Value = -21088
$ws.Cells.Item(26, 13).ClearContents()

$ws.Cells.Item(35, 8).Value = 11000
$ws.Cells.Item(35, 10).Value = 14266.667
$ws.Cells.Item(35, 12).Value = 14266.667
$ws.Cells.Item(35, 14).Value = -15024.667

$ws.Cells.Item(86, 8).Value = 7355500
$ws.Cells.Item(86, 9).Value = 20834934
$ws.Cells.Item(86, 10).Value = 3081.818
$ws.Cells.Item(86, 11).Value = 20834934
$ws.Cells.Item(86, 12).Value = 3081.818
$ws.Cells.Item(86, 13).Value = -20833811
$ws.Cells.Item(86, 14).Value = -5327.818

$ws.Cells.Item(89, 8).Value = 7355500
$ws.Cells.Item(89, 9).Value = 20834934
$ws.Cells.Item(89, 10).Value = 3081.818
$ws.Cells.Item(89, 11).Value = 104174670
$ws.Cells.Item(89, 12).Value = 15409.09
$ws.Cells.Item(89, 13).Value = -104169054
$ws.Cells.Item(89, 14).Value = -26641.09

$ws.Cells.Item(111, 8).Value = 7118.4287
$ws.Cells.Item(111, 9).Value = 9205.799999999999
$ws.Cells.Item(111, 10).Value = 1900
$ws.Cells.Item(111, 11).Value = 27617.4
$ws.Cells.Item(111, 12).Value = 5700
$ws.Cells.Item(111, 13).Value = -24550.4
$ws.Cells.Item(111, 14).Value = -11834

$ws.Cells.Item(132, 8).Value = 1829.1082
$ws.Cells.Item(132, 9).Value = 1911.0883
$ws.Cells.Item(132, 10).Value = 900
$ws.Cells.Item(132, 11).Value = 5733.2649
$ws.Cells.Item(132, 12).Value = 2700
$ws.Cells.Item(132, 13).Value = -3203.2649
$ws.Cells.Item(132, 14).Value = -7760

$ws.Cells.Item(135, 8).Value = 75001220
$ws.Cells.Item(135, 9).Value = 27779072
$ws.Cells.Item(135, 10).Value = 500000500
$ws.Cells.Item(135, 11).Value = 250011648
$ws.Cells.Item(135, 12).Value = 4500004500
$ws.Cells.Item(135, 13).Value = -250009113
$ws.Cells.Item(135, 14).Value = -4500009570

$ws.Cells.Item(137, 8).Value = 2965.238
$ws.Cells.Item(137, 9).Value = 3060.6365
$ws.Cells.Item(137, 10).Value = 2860.3
$ws.Cells.Item(137, 11).Value = 9181.9095
$ws.Cells.Item(137, 12).Value = 8580.900000000001
$ws.Cells.Item(137, 13).Value = -6631.9095
$ws.Cells.Item(137, 14).Value = -13680.9

$ws.Cells.Item(138, 8).Value = 2759.8857
$ws.Cells.Item(138, 9).Value = 1121.909
$ws.Cells.Item(138, 10).Value = 5531.846
$ws.Cells.Item(138, 11).Value = 3365.727
$ws.Cells.Item(138, 12).Value = 16595.538
$ws.Cells.Item(138, 13).Value = 1774.273
$ws.Cells.Item(138, 14).Value = -26875.538

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(26, 8).Value = 1700
$ws.Cells.Item(26, 9).Value = 1700
$ws.Cells.Item(26, 11).Value = 1700
$ws.Cells.Item(26, 13).Value = -1370

$ws.Cells.Item(32, 8).Value = 30401.303
$ws.Cells.Item(32, 9).Value = 34515.97
$ws.Cells.Item(32, 10).Value = 14857
$ws.Cells.Item(32, 11).Value = 34515.97
$ws.Cells.Item(32, 12).Value = 14857
$ws.Cells.Item(32, 13).Value = -34228.97
$ws.Cells.Item(32, 14).Value = -15431

$ws.Cells.Item(61, 8).Value = 6738.0195
$ws.Cells.Item(61, 9).Value = 3331.1025
$ws.Cells.Item(61, 10).Value = 17810.5
$ws.Cells.Item(61, 11).Value = 3331.1025
$ws.Cells.Item(61, 12).Value = 17810.5
$ws.Cells.Item(61, 13).Value = -3119.1025
$ws.Cells.Item(61, 14).Value = -18234.5

$ws.Cells.Item(74, 8).Value = 5173.2583
$ws.Cells.Item(74, 9).Value = 2250.05
$ws.Cells.Item(74, 10).Value = 10488.182
$ws.Cells.Item(74, 11).Value = 2250.05
$ws.Cells.Item(74, 12).Value = 10488.182
$ws.Cells.Item(74, 13).Value = -1376.05
$ws.Cells.Item(74, 14).Value = -12236.182

$ws.Cells.Item(77, 8).Value = 5173.2583
$ws.Cells.Item(77, 9).Value = 2250.05
$ws.Cells.Item(77, 10).Value = 10488.182
$ws.Cells.Item(77, 11).Value = 11250.25
$ws.Cells.Item(77, 12).Value = 52440.91
$ws.Cells.Item(77, 13).Value = -6882.25
$ws.Cells.Item(77, 14).Value = -61176.91

$ws.Cells.Item(102, 8).Value = 1124861.4
$ws.Cells.Item(102, 9).Value = 1685587
$ws.Cells.Item(102, 10).Value = 3410.0908
$ws.Cells.Item(102, 11).Value = 1685587
$ws.Cells.Item(102, 12).Value = 3410.0908
$ws.Cells.Item(102, 13).Value = -1683965
$ws.Cells.Item(102, 14).Value = -6654.0908

$ws.Cells.Item(122, 8).Value = 2604.5334
$ws.Cells.Item(122, 9).Value = 2926.6667
$ws.Cells.Item(122, 10).Value = 2389.7778
$ws.Cells.Item(122, 11).Value = 8780.000100000001
$ws.Cells.Item(122, 12).Value = 7169.3334
$ws.Cells.Item(122, 13).Value = -6330.000100000001
$ws.Cells.Item(122, 14).Value = -12069.3334

$ws.Cells.Item(128, 8).Value = 57500
$ws.Cells.Item(128, 10).Value = 57500
$ws.Cells.Item(128, 12).Value = 57500
$ws.Cells.Item(128, 14).Value = -67460

$ws.Cells.Item(136, 8).Value = 6738.0195
$ws.Cells.Item(136, 9).Value = 3331.1025
$ws.Cells.Item(136, 10).Value = 17810.5
$ws.Cells.Item(136, 11).Value = 9993.307499999999
$ws.Cells.Item(136, 12).Value = 53431.5
$ws.Cells.Item(136, 13).Value = -7443.307499999999
$ws.Cells.Item(136, 14).Value = -58531.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1869.6666
$ws.Cells.Item(20, 9).Value = 1600
$ws.Cells.Item(20, 10).Value = 2409
$ws.Cells.Item(20, 11).Value = 1600
$ws.Cells.Item(20, 12).Value = 2409
$ws.Cells.Item(20, 13).Value = -1353
$ws.Cells.Item(20, 14).Value = -2903

$ws.Cells.Item(62, 8).Value = 36181
$ws.Cells.Item(62, 10).Value = 36181
$ws.Cells.Item(62, 12).Value = 36181
$ws.Cells.Item(62, 14).Value = -37553

$ws.Cells.Item(65, 8).Value = 36181
$ws.Cells.Item(65, 10).Value = 36181
$ws.Cells.Item(65, 12).Value = 108543
$ws.Cells.Item(65, 14).Value = -115407

$ws.Cells.Item(99, 8).Value = 972.8570999999999
$ws.Cells.Item(99, 9).Value = 880
$ws.Cells.Item(99, 11).Value = 880
$ws.Cells.Item(99, 13).Value = 618

$ws.Cells.Item(134, 8).Value = 40302.777
$ws.Cells.Item(134, 9).Value = 3381.1667
$ws.Cells.Item(134, 10).Value = 114146
$ws.Cells.Item(134, 11).Value = 10143.5001
$ws.Cells.Item(134, 12).Value = 342438
$ws.Cells.Item(134, 13).Value = -7608.500100000001
$ws.Cells.Item(134, 14).Value = -347508

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4746.9023
$ws.Cells.Item(31, 9).Value = 5116.643
$ws.Cells.Item(31, 10).Value = 3950.5386
$ws.Cells.Item(31, 11).Value = 5116.643
$ws.Cells.Item(31, 12).Value = 3950.5386
$ws.Cells.Item(31, 13).Value = -4821.643
$ws.Cells.Item(31, 14).Value = -4540.5386

$ws.Cells.Item(34, 8).Value = 4746.9023
$ws.Cells.Item(34, 9).Value = 5116.643
$ws.Cells.Item(34, 10).Value = 3950.5386
$ws.Cells.Item(34, 11).Value = 5116.643
$ws.Cells.Item(34, 12).Value = 3950.5386
$ws.Cells.Item(34, 13).Value = -4914.643
$ws.Cells.Item(34, 14).Value = -4354.5386

$ws.Cells.Item(58, 8).Value = 2275672.5
$ws.Cells.Item(58, 9).Value = 4134537.5
$ws.Cells.Item(58, 10).Value = 3726.111
$ws.Cells.Item(58, 11).Value = 4134537.5
$ws.Cells.Item(58, 12).Value = 3726.111
$ws.Cells.Item(58, 13).Value = -4134334.5
$ws.Cells.Item(58, 14).Value = -4132.111

$ws.Cells.Item(134, 8).Value = 2588.6948
$ws.Cells.Item(134, 9).Value = 1729.7
$ws.Cells.Item(134, 10).Value = 3477.3103
$ws.Cells.Item(134, 11).Value = 5189.1
$ws.Cells.Item(134, 12).Value = 10431.9309
$ws.Cells.Item(134, 13).Value = -2654.1
$ws.Cells.Item(134, 14).Value = -15501.9309

$ws.Cells.Item(136, 8).Value = 2275672.5
$ws.Cells.Item(136, 9).Value = 4134537.5
$ws.Cells.Item(136, 10).Value = 3726.111
$ws.Cells.Item(136, 11).Value = 12403612.5
$ws.Cells.Item(136, 12).Value = 11178.333
$ws.Cells.Item(136, 13).Value = -12401062.5
$ws.Cells.Item(136, 14).Value = -16278.333

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(6, 8).Value = 117
$ws.Cells.Item(6, 9).Value = 117
$ws.Cells.Item(6, 10).Value = 0
$ws.Cells.Item(6, 11).Value = 351
$ws.Cells.Item(6, 12).Value = 0
$ws.Cells.Item(6, 13).Value = -238
$ws.Cells.Item(6, 14).ClearContents()

$ws.Cells.Item(104, 8).Value = 2229
$ws.Cells.Item(104, 10).Value = 2229
$ws.Cells.Item(104, 12).Value = 6687
$ws.Cells.Item(104, 14).Value = -11929

$ws.Cells.Item(117, 8).Value = 0
$ws.Cells.Item(117, 9).Value = 0
$ws.Cells.Item(117, 10).Value = 0
$ws.Cells.Item(117, 11).Value = 0
$ws.Cells.Item(117, 12).Value = 0
$ws.Cells.Item(117, 13).ClearContents()
$ws.Cells.Item(117, 14).ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(43, 8).Value = 2017
$ws.Cells.Item(43, 10).Value = 0
$ws.Cells.Item(43, 12).Value = 0
$ws.Cells.Item(43, 14).ClearContents()

$ws.Cells.Item(132, 8).Value = 8549.406000000001
$ws.Cells.Item(132, 9).Value = 5817.7407
$ws.Cells.Item(132, 11).Value = 17453.2221
$ws.Cells.Item(132, 13).Value = -14923.2221

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 931.9524
$ws.Cells.Item(16, 9).Value = 627.1667
$ws.Cells.Item(16, 10).Value = 2760.6667
$ws.Cells.Item(16, 11).Value = 627.1667
$ws.Cells.Item(16, 12).Value = 2760.6667
$ws.Cells.Item(16, 13).Value = -457.1667
$ws.Cells.Item(16, 14).Value = -3100.6667

$ws.Cells.Item(59, 8).Value = 0
$ws.Cells.Item(59, 10).Value = 0
$ws.Cells.Item(59, 12).Value = 0
$ws.Cells.Item(59, 14).ClearContents()

$ws.Cells.Item(82, 8).Value = 2546.0588
$ws.Cells.Item(82, 9).Value = 2465.4167
$ws.Cells.Item(82, 11).Value = 2465.4167
$ws.Cells.Item(82, 13).Value = -2104.4167

$ws.Cells.Item(85, 8).Value = 2546.0588
$ws.Cells.Item(85, 9).Value = 2465.4167
$ws.Cells.Item(85, 11).Value = 2465.4167
$ws.Cells.Item(85, 13).Value = -1217.4167

$ws.Cells.Item(132, 8).Value = 6723.884
$ws.Cells.Item(132, 9).Value = 8284.259
$ws.Cells.Item(132, 11).Value = 24852.777
$ws.Cells.Item(132, 13).Value = -22322.777

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(54, 8).Value = 14880.333
$ws.Cells.Item(54, 10).Value = 14880.333
$ws.Cells.Item(54, 12).Value = 14880.333
$ws.Cells.Item(54, 14).Value = -15920.333

$ws.Cells.Item(122, 8).Value = 2445.7273
$ws.Cells.Item(122, 9).Value = 1589.8182
$ws.Cells.Item(122, 10).Value = 3301.6365
$ws.Cells.Item(122, 11).Value = 4769.4546
$ws.Cells.Item(122, 12).Value = 9904.9095
$ws.Cells.Item(122, 13).Value = -2319.4546
$ws.Cells.Item(122, 14).Value = -14804.9095

$ws.Cells.Item(126, 8).Value = 1504.3871
$ws.Cells.Item(126, 9).Value = 1746
$ws.Cells.Item(126, 10).Value = 1246.6666
$ws.Cells.Item(126, 11).Value = 5238
$ws.Cells.Item(126, 12).Value = 3739.9998
$ws.Cells.Item(126, 13).Value = -2768
$ws.Cells.Item(126, 14).Value = -8679.9998

$ws.Cells.Item(132, 8).Value = 2754.5186
$ws.Cells.Item(132, 9).Value = 1497.6364
$ws.Cells.Item(132, 10).Value = 3618.625
$ws.Cells.Item(132, 11).Value = 4492.9092
$ws.Cells.Item(132, 12).Value = 10855.875
$ws.Cells.Item(132, 13).Value = -1962.9092
$ws.Cells.Item(132, 14).Value = -15915.875

$ws.Cells.Item(136, 8).Value = 6925.787
$ws.Cells.Item(136, 9).Value = 7152.2
$ws.Cells.Item(136, 10).Value = 6758.074
$ws.Cells.Item(136, 11).Value = 21456.6
$ws.Cells.Item(136, 12).Value = 20274.222
$ws.Cells.Item(136, 13).Value = -20274.222
$ws.Cells.Item(136, 14).Value = -25374.222
